# faturamento_diario.xlsx – reload Sheet1 with the refreshed ADD daily-revenue export
# (modelagem RFMA / recorrencia / retencao / faturamento migrada para o banco padrao;
# dados novos de ADD carregados, por isso os valores de B/dia mudam e a planilha cresce
# de 62 para 65 linhas de dados).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dia, total_venda, Mes, Ano, Periodo per row (rows 2..66)
$rows = @(
  @(2, 20078.49, 5, 2025, "05/2025"),
  @(5, 23304.62, 5, 2025, "05/2025"),
  @(6, 18561.83, 5, 2025, "05/2025"),
  @(7, 21063.46, 5, 2025, "05/2025"),
  @(8, 19773.95, 5, 2025, "05/2025"),
  @(9, 32201.9, 5, 2025, "05/2025"),
  @(1, 29887.28, 4, 2025, "04/2025"),
  @(2, 48708.66, 4, 2025, "04/2025"),
  @(3, 35472.5, 4, 2025, "04/2025"),
  @(4, 7001.320000000001, 4, 2025, "04/2025"),
  @(7, 7753.57, 4, 2025, "04/2025"),
  @(8, 7630.52, 4, 2025, "04/2025"),
  @(9, 19608.04, 4, 2025, "04/2025"),
  @(10, 18202.05, 4, 2025, "04/2025"),
  @(11, 9551.92, 4, 2025, "04/2025"),
  @(14, 34424.2, 4, 2025, "04/2025"),
  @(15, 12545.33, 4, 2025, "04/2025"),
  @(16, 15403.61, 4, 2025, "04/2025"),
  @(17, 13921.35, 4, 2025, "04/2025"),
  @(22, 16758.86, 4, 2025, "04/2025"),
  @(23, 21276.73, 4, 2025, "04/2025"),
  @(24, 17634.6, 4, 2025, "04/2025"),
  @(25, 6381.17, 4, 2025, "04/2025"),
  @(28, 30330.68, 4, 2025, "04/2025"),
  @(29, 31647.45, 4, 2025, "04/2025"),
  @(30, 31765.2, 4, 2025, "04/2025"),
  @(3, 19552.17, 3, 2025, "03/2025"),
  @(6, 35882.61, 3, 2025, "03/2025"),
  @(7, 13042.5, 3, 2025, "03/2025"),
  @(10, 19948.65, 3, 2025, "03/2025"),
  @(11, 8539.889999999999, 3, 2025, "03/2025"),
  @(12, 25156.89, 3, 2025, "03/2025"),
  @(13, 14118.35, 3, 2025, "03/2025"),
  @(14, 19961.1, 3, 2025, "03/2025"),
  @(17, 24428.36, 3, 2025, "03/2025"),
  @(18, 3443.92, 3, 2025, "03/2025"),
  @(19, 17409.2, 3, 2025, "03/2025"),
  @(20, 4729.55, 3, 2025, "03/2025"),
  @(21, 22027.47, 3, 2025, "03/2025"),
  @(24, 26877.71, 3, 2025, "03/2025"),
  @(25, 24190.81, 3, 2025, "03/2025"),
  @(26, 18815.44, 3, 2025, "03/2025"),
  @(27, 13162.5, 3, 2025, "03/2025"),
  @(28, 9765.5, 3, 2025, "03/2025"),
  @(31, 41909.04, 3, 2025, "03/2025"),
  @(3, 26998.77, 2, 2025, "02/2025"),
  @(4, 18480.02, 2, 2025, "02/2025"),
  @(5, 21525.35, 2, 2025, "02/2025"),
  @(6, 29313.91, 2, 2025, "02/2025"),
  @(7, 20990.36, 2, 2025, "02/2025"),
  @(10, 26256.4, 2, 2025, "02/2025"),
  @(11, 9265.960000000001, 2, 2025, "02/2025"),
  @(12, 22580.88, 2, 2025, "02/2025"),
  @(13, 9109.67, 2, 2025, "02/2025"),
  @(14, 7871.8, 2, 2025, "02/2025"),
  @(17, 8339.190000000001, 2, 2025, "02/2025"),
  @(18, 17666.8, 2, 2025, "02/2025"),
  @(19, 19478, 2, 2025, "02/2025"),
  @(20, 7205.39, 2, 2025, "02/2025"),
  @(21, 8843.360000000001, 2, 2025, "02/2025"),
  @(24, 20183.37, 2, 2025, "02/2025"),
  @(25, 12027.74, 2, 2025, "02/2025"),
  @(26, 53797.58, 2, 2025, "02/2025"),
  @(27, 39559.92, 2, 2025, "02/2025"),
  @(28, 77821.24000000001, 2, 2025, "02/2025")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  $ws.Range("A$r").Value = $row[0]
  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
  $ws.Range("E$r").Value = $row[4]
}

"Updated $($rows.Count) data rows (2..$($rows.Count + 1)) on Sheet1"
